$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix the capitalization of the "yami" entry (row 8, column A)
$ws.Range("A8").Value = "Yami"

# Mark TipoPrestacion (D7) and FacturaCliente (E7) as implemented classes,
# matching the green font + border already used for B7/C7
$ws.Range("D7:E7").Font.Color = 5287936

# Move the active selection to A10
$ws.Range("A10").Select() | Out-Null
